# Update "yearly" overview sheet:
#  - Roll the 5-year reporting window forward by one year:
#       1396,1397,1398,1399,1400  ->  1397,1398,1399,1400,1401
#    (drop the oldest column's data/label, shift the rest left, append new 1401 data)
#  - This is the "update database" part of the commit; the header labels in
#    columns E:I on rows 8 and 24 shift accordingly and a brand-new 1401
#    figure is added as the rightmost (I) column for every metric row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column headers (row 8 and row 24), shifted one year forward ---
$headers = @(
    "دوازده ماهه منتهی به 1397/12",
    "دوازده ماهه منتهی به 1398/12",
    "دوازده ماهه منتهی به 1399/12",
    "دوازده ماهه منتهی به 1400/12",
    "دوازده ماهه منتهی به 1401/12"
)

for ($i = 0; $i -lt 5; $i++) {
    $col = 5 + $i   # E=5 .. I=9
    $ws.Cells.Item(8, $col).Value = $headers[$i]
    $ws.Cells.Item(24, $col).Value = $headers[$i]
}

# --- Data rows: each row's old F:I values shift left into E:H
#     (old F->E, G->F, H->G, I->H) and a brand new figure lands in I. ---
$rowData = @{
    10 = @(12146, 92812, 198138, 278983, 567750)
    11 = @(0, 0, 0, 0, 0)
    12 = @(0, 0, 0, 0, 0)
    13 = @(421749, 160631, 157987, 485445, 1081871)
    14 = @(0, 0, 0, 0, 0)
    15 = @(0, 0, 0, 0, 0)
    16 = @(3981, 9138, 19816, 16748, 28450)
    17 = @(263058, 652426, 1061184, 1898964, 2445425)
    18 = @(0, 0, 0, 0, 0)
    19 = @(99148, 181673, 353195, 404487, 641826)
    20 = @(800082, 1096680, 1790320, 3084627, 4765322)
    26 = @(794, 856, 1190, 1235, 1300)
    27 = @(264, 320, 578, 571, 600)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt 5; $i++) {
        $col = 5 + $i   # E=5 .. I=9
        $ws.Cells.Item($row, $col).Value = $vals[$i]
    }
}
